$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the inventory values: web_servers host + fix db_servers row
$ws.Range("B2").Value = "web1"
$ws.Range("A3").Value = "db_servers"
$ws.Range("B3").Value = "db1"

# Restore selection to B3 and set the sheet zoom, matching the saved view state
$ws.Range("B3").Select()
$excel.ActiveWindow.Zoom = 136
